# Generate Report for Archive
#
# 1. Every cell whose status was "Ready for handoff" moves to "In Translation"
#    (Overview!E2:F2/E3:F3, zh-cn!C2:C3, de-de!C2:C3 - all eight cells share
#    the same underlying string).
# 2. The zh-cn / de-de status columns (and the matching columns on the
#    Overview summary sheet) are narrowed to their new autofit width.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# Narrower column width (was ~17.22 chars stored width, now ~13.41).
# ColumnWidth is expressed in characters and Excel snaps it to whole
# pixels, so 12.5 is the input that lands on the closest achievable
# stored width to the target.
$newColumnWidth = 12.5

$wsOverview.Range("E1").ColumnWidth = $newColumnWidth
$wsOverview.Range("F1").ColumnWidth = $newColumnWidth

$wsZhCn.Range("C1").ColumnWidth = $newColumnWidth

$wsDeDe.Range("C1").ColumnWidth = $newColumnWidth
